$d = $word.ActiveDocument

# 1. Insert "self-hosted in Azure with GitHub CI/CD, " right before
#    "for an interactive resume with more information, and " so the
#    sentence reads "..., self-hosted in Azure with GitHub CI/CD, for an
#    interactive resume with more information, and ...".
#    We use a collapsed (zero-length) insertion-point range rather than
#    replacing the whole phrase so the untouched surrounding text keeps
#    its original (un-styled) run formatting instead of inheriting the
#    formatting of whatever neighboring run a plain Replace would pick up.
$findRange = $d.Content
$findRange.Find.ClearFormatting()
$findRange.Find.Text = "for an interactive resume with more information, and "
$found = $findRange.Find.Execute()
if ($found) {
    $cursor = $findRange.Duplicate
    $cursor.Collapse(1)
    $cursor.InsertBefore("self-hosted in Azure with GitHub CI/CD, ")
    Write-Host "Inserted 'self-hosted in Azure with GitHub CI/CD, ' before the interactive-resume sentence."
} else {
    Write-Host "WARNING: could not find the interactive-resume sentence."
}

# 2. Replace the "Git" skill bullet with "Docker".
$gitRange = $d.Content
$gitRange.Find.ClearFormatting()
$gitRange.Find.Text = "Git"
$gitRange.Find.MatchWholeWord = $true
$gitRange.Find.MatchCase = $true
$foundGit = $gitRange.Find.Execute()
if ($foundGit) {
    $gitRange.Text = "Docker"
    Write-Host "Replaced the 'Git' skill bullet with 'Docker'."
} else {
    Write-Host "WARNING: could not find the 'Git' skill bullet."
}
